$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new cells P1, Q1 (styled like the rest of row 1, style index 1 => bold/center/border)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Row 2: new cells P2, Q2 (no special style)
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0

# Rows 3-25: swap I/K and M/O columns' values, and add new P,Q columns with value 2
for ($r = 3; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
